$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing D:K data to F:M.
$ws.Range("D:E").Insert(-4161)

# Carry over number formatting from the (now-shifted) original columns
# so the new D/E columns match the date row (style 2) and data rows (style 3).
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("G7:G102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newData = @(
    @{Row=7; D=43435; E=43344},
    @{Row=8; D=413200; E=355300},
    @{Row=9; D=282700; E=261800},
    @{Row=10; D=130500; E=93500},
    @{Row=11; D=$null; E=$null},
    @{Row=12; D="NA"; E="NA"},
    @{Row=13; D=0; E=0},
    @{Row=14; D=0; E=0},
    @{Row=15; D=12400; E=12800},
    @{Row=16; D=$null; E=$null},
    @{Row=17; D=442200; E=417800},
    @{Row=18; D=-29000; E=-62500},
    @{Row=19; D=$null; E=$null},
    @{Row=20; D=-100; E=900},
    @{Row=21; D=-14700; E=-46800},
    @{Row=22; D=3500; E=3600},
    @{Row=23; D=-32600; E=-65200},
    @{Row=24; D=17900; E=-14100},
    @{Row=25; D=0; E=0},
    @{Row=26; D=-50400; E=-51100},
    @{Row=27; D=-50400; E=-51100},
    @{Row=28; D=0; E=0},
    @{Row=29; D="NA"; E="NA"},
    @{Row=30; D=0; E=0},
    @{Row=31; D=0; E=0},
    @{Row=32; D=100; E=-900},
    @{Row=33; D=-50400; E=-51100},
    @{Row=34; D=0; E=0},
    @{Row=35; D=-50400; E=-51100},
    @{Row=38; D=43435; E=43344},
    @{Row=39; D=$null; E=$null},
    @{Row=40; D=$null; E=$null},
    @{Row=41; D=7800; E=27600},
    @{Row=42; D=63300; E=89200},
    @{Row=43; D=36300; E=24200},
    @{Row=44; D=388300; E=386700},
    @{Row=45; D=56700; E=51800},
    @{Row=46; D=552300; E=579400},
    @{Row=47; D=0; E=0},
    @{Row=48; D=159700; E=168100},
    @{Row=49; D=0; E=0},
    @{Row=50; D=0; E=0},
    @{Row=51; D=0; E=0},
    @{Row=52; D=33300; E=57500},
    @{Row=53; D=0; E=0},
    @{Row=54; D=745300; E=805000},
    @{Row=55; D=$null; E=$null},
    @{Row=56; D=$null; E=$null},
    @{Row=57; D=172700; E=181500},
    @{Row=58; D=2000; E=2000},
    @{Row=59; D=162300; E=161100},
    @{Row=60; D=337000; E=344600},
    @{Row=61; D=197000; E=197300},
    @{Row=62; D=54100; E=55900},
    @{Row=63; D=0; E=0},
    @{Row=64; D=0; E=0},
    @{Row=65; D=0; E=0},
    @{Row=66; D=588100; E=597800},
    @{Row=67; D=$null; E=$null},
    @{Row=68; D=0; E=0},
    @{Row=69; D=0; E=0},
    @{Row=70; D=0; E=0},
    @{Row=71; D=0; E=0},
    @{Row=72; D=603200; E=653700},
    @{Row=73; D=0; E=0},
    @{Row=74; D=0; E=0},
    @{Row=75; D=0; E=0},
    @{Row=76; D=157200; E=207200},
    @{Row=77; D=0; E=0},
    @{Row=80; D=43435; E=43344},
    @{Row=81; D=-50400; E=-51100},
    @{Row=82; D=$null; E=$null},
    @{Row=83; D=14400; E=14900},
    @{Row=84; D=0; E=0},
    @{Row=85; D=0; E=0},
    @{Row=86; D=0; E=0},
    @{Row=87; D=0; E=0},
    @{Row=88; D=0; E=0},
    @{Row=89; D=-42900; E=-28700},
    @{Row=90; D=$null; E=$null},
    @{Row=91; D=-5800; E=-13500},
    @{Row=92; D=0; E=0},
    @{Row=93; D=0; E=0},
    @{Row=94; D=-2400; E=-11200},
    @{Row=95; D=$null; E=$null},
    @{Row=96; D=0; E=0},
    @{Row=97; D=0; E=0},
    @{Row=98; D=0; E=0},
    @{Row=99; D=0; E=0},
    @{Row=100; D=-100; E=-100},
    @{Row=101; D=-300; E=0},
    @{Row=102; D=-45700; E=-40000}
)

foreach ($item in $newData) {
    $r = $item.Row
    if ($null -ne $item.D) {
        $ws.Cells.Item($r, 4).Value = $item.D
    }
    if ($null -ne $item.E) {
        $ws.Cells.Item($r, 5).Value = $item.E
    }
}
